$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the degree/certification abbreviations to include periods
$ws.Range("A15").Value = "Ph.D. Landscape Ecology"
$ws.Range("A16").Value = "M.Sc. Behavioural Ecology"
$ws.Range("A17").Value = "B.Ed. Education"
$ws.Range("A18").Value = "B.Sc."

# Reflect the active cell selection change seen in the saved file
$ws.Range("A18").Select()
